$d = $word.ActiveDocument

# The paragraph "Not only eo option in postproc vols, also lp" is removed in
# its entirety - its text plus its trailing paragraph mark - which merges the
# following paragraph ("Exception when partial stacks of wrong ldim") into
# it, exactly like selecting the whole line (including the pilcrow) in Word
# and pressing Delete.

$targetText = "Not only eo option in postproc vols, also lp"

$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    $t = $t.TrimEnd([char]13, [char]7)
    if ($t -eq $targetText) {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $start = $target.Range.Start
    $nextPara = $target.Next()
    $end = $nextPara.Range.Start

    $r = $d.Range($start, $end)
    $r.Delete()

    # Word keeps the "_GoBack" bookmark at the location of the last edit;
    # relocate it here (it previously sat at the very end of the list, right
    # after the "mailx ..." text).
    $bmRange = $d.Range($start, $start)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
